$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $found = $d.Content.Find.Execute(
        $findText,
        $false, $false, $false, $false, $false, $true, 1, $false,
        $replaceText,
        2)
    if (-not $found) {
        Write-Output ("WARNING: text not found -> " + $findText)
    }
}

# 1. "Uno dei motivi potrebbe essere..." -- expand "degli imprenditori e le loro situazioni di
#    vita." to "degli imprenditori,  le loro situazioni di vita e la loro situazione sociale."
Replace-Text `
    "date le competenze degli imprenditori e le loro situazioni di vita." `
    "date le competenze degli imprenditori,  le loro situazioni di vita e la loro situazione sociale."

# 2. "principalmente per 3 motivi:" -> "principalmente per tre motivi:"
Replace-Text "per 3 motivi:" "per tre motivi:"

# 3. "ci ritrova in una situazione peggiore" -> "ci si ritrova in una situazione peggiore"
Replace-Text `
    "concesso, ci ritrova in una situazione peggiore rispetto a" `
    "concesso, ci si ritrova in una situazione peggiore rispetto a"

# 4. "quella iniziale, cioè ci si trova a dover rimborsare" -> "quella iniziale, cioè si deve rimborsare"
Replace-Text `
    "quella iniziale, cioè ci si trova a dover rimborsare il debito maggiorato di interessi (anche se bassi)." `
    "quella iniziale, cioè si deve rimborsare il debito maggiorato di interessi (anche se bassi)."

# 5. Merge "In conclusione" + " pensiamo che..." into a single run without the proofErr wrapper.
Replace-Text `
    "In conclusione pensiamo che, finché sussistono questi problemi strutturali, il benessere che la concessione" `
    "In conclusione pensiamo che, finché sussistono questi problemi strutturali, il benessere che la concessione"

# 6. Remove the trailing ")" run that closed the parenthetical note at the end of the document.
Replace-Text `
    "seguendo soprattutto i nuovi imprenditori.)" `
    "seguendo soprattutto i nuovi imprenditori."
